# "updated data mid-morning 7/27" — append new Activity-log rows (391-410)
# to the "2019" sheet / Table2, covering 7/24-7/27/2019.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")
$tbl = $ws.ListObjects.Item(1)

$formula = '=IF(Table2[[#This Row],[Activity]]="Sleep",(Table2[[#This Row],[End]]-Table2[[#This Row],[Start]])*24,"NA")'

# Copy number-format only (no value/content) from a template cell so the new
# cell reuses an existing style entry instead of minting a fresh one.
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)
}

# Add one data row. $start/$end are Excel serial date-times ($end can be $null
# for rows with no End value). $activity is the Activity column text,
# $comment is the Comment column text (or $null when absent).
function Add-ActivityRow {
    param(
        [int]$Row,
        [double]$Start,
        $End,
        [string]$Activity,
        $Comment,
        [string]$StyleTemplateAddr
    )

    $colA = "A{0}" -f $Row
    Copy-Format $StyleTemplateAddr $colA
    $ws.Cells.Item($Row, 1).Value = $Start

    if ($null -ne $End) {
        $colB = "B{0}" -f $Row
        Copy-Format $StyleTemplateAddr $colB
        $ws.Cells.Item($Row, 2).Value = $End
    }

    $ws.Cells.Item($Row, 3).Value = $Activity

    if ($null -ne $Comment) {
        $ws.Cells.Item($Row, 4).Value = $Comment
    }

    $ws.Cells.Item($Row, 5).Formula = $formula
}

# Rows 391-410, transcribed from the source workbook's new entries.
# (Positional args: Row, Start, End, Activity, Comment, StyleTemplateAddr.)
Add-ActivityRow 391 43670.357638888891 $null "Food" "English muffin w/egg strawberry" "A390"
Add-ActivityRow 392 43670.891435185185 43671.227083333331 "Sleep" $null "A390"
Add-ActivityRow 393 43670.684027777781 $null "Food" "Blueberry almond smoothie " "A390"
Add-ActivityRow 394 43670.78125 $null "Food" "Cod + bread + wine" "A390"
Add-ActivityRow 395 43670.635416666664 43670.666666666664 "Exercise" $null "A390"
Add-ActivityRow 396 43670.53125 43670.5625 "Exercise" $null "A390"
Add-ActivityRow 397 43670.5 $null "Food" "Pad Thai + quinoa" "A390"
Add-ActivityRow 398 43671.920648148145 43672.236805555556 "Sleep" $null "A390"
Add-ActivityRow 399 43671.291666666664 $null "Food" "Bread + egg" "A390"
Add-ActivityRow 400 43671.322916666664 43671.354166666664 "Exercise" $null "A390"
Add-ActivityRow 401 43671.354166666664 $null "Food" "Smoothie" "A390"
Add-ActivityRow 402 43671.520833333336 $null "Food" "Beans + salad + seafood + bread" "A390"
Add-ActivityRow 403 43671.791666666664 $null "Food" "Seafood + rice + salad + brownie" "A390"
Add-ActivityRow 404 43671.833333333336 43671.868055555555 "Exercise" $null "A390"
Add-ActivityRow 405 43672.364583333336 $null "Food" "Bread + egg" "A390"
Add-ActivityRow 406 43672.770833333336 $null "Food" "Salad w/pork belly and bread" "A390"
Add-ActivityRow 407 43672.510416666664 $null "Food" "Cabbage + tofu" "A390"
Add-ActivityRow 408 43672.909722222219 43673.236111111109 "Sleep" $null "A252"
Add-ActivityRow 409 43672.409722222219 43672.420138888891 "Exercise" $null "A390"
Add-ActivityRow 410 43673.28125 $null "Food" "English muffin" "A390"

# Grow the table + worksheet dimension to cover the new rows.
$tbl.Resize($ws.Range("A1:E410"))

# Match the author's final viewport/selection state.
$ws.Range("A411").Select()
$excel.ActiveWindow.ScrollRow = 379

Write-Host "Added rows 391-410 to Table2 on sheet '2019'."
